$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 (Герюгов Ислам): fill in ДЗ1-ДЗ6 with 5, add ДЗ7 (I8) with matching
#     border/alignment format copied from an existing "has ДЗ7" row (I6) ---
$ws.Range("C8:H8").Value = 5
$ws.Range("I6").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("I8").Value = 5

# --- Row 22 (Рустан Вячеслав): fill in ДЗ1-ДЗ3 with 5 ---
$ws.Range("C22:E22").Value = 5

# --- Row 24 (Сидаков Амир): fill in ДЗ1-ДЗ2 with 5 ---
$ws.Range("C24:D24").Value = 5

# --- Row 26 (Теплюк Дмитрий): fill in ДЗ1-ДЗ6 with 5, add ДЗ7 (I26) with
#     matching border/alignment format copied from I13 ---
$ws.Range("C26:H26").Value = 5
$ws.Range("I13").Copy()
$ws.Range("I26").PasteSpecial(-4122)
$ws.Range("I26").Value = 5

# --- Restore the active selection/view state: the frozen pane's top-left
#     cell and the active cell in the bottom-right pane ---
$ws.Range("I8").Select()

$excel.CutCopyMode = 0
